# Rename the PP-impeller value-input headers on the first sheet
# ("ecological_params") to the shorter, consistent naming used
# elsewhere in the workbook, drop the cell border that used to set
# them apart (matching the "Scaling" sheet's header formatting), and
# restore the view/selection state left behind from editing
# ("formatting of value input files").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename headers B1:D1 on the "ecological_params" sheet.
$ws1.Range("B1").Value = "PP_virgin"
$ws1.Range("C1").Value = "PP_recycled"
$ws1.Range("D1").Value = "PP_recycled_industrial"

# Drop the border on those header cells so they match the plain
# (borderless) header style used on the "Scaling" sheet.
$ws1.Range("B1:D1").Borders.LineStyle = 0

# Leave the workbook with "ecological_params" as the active sheet and
# cell D11 selected there (mirrors the selection state on "Scaling").
$ws1.Activate()
$ws1.Range("D11").Select()
